# Anti-foam job aid.docx - SharePoint content-type / managed-metadata resync
#
# The document was reconnected to a different SharePoint content type
# ("Document" 0x0101004A0300... -> 0x010100FE4CDFAD...). Word rewrites the
# hidden custom XML parts that carry the content-type schema (customXml
# item2.xml), the document-management property values (item3.xml) and the
# content-type/taxonomy sync marker (item4.xml) to match the new site
# columns (Content Author Email, Content Owner, Record, Sensitivity,
# Language, Group, Project Name, Tags, ...).
#
# This script reproduces that resync using the Word object model:
# the old custom XML parts (content-type schema + document management
# properties, under the legacy "700a75ce-..." / "ce04ccd7-..." namespaces)
# are removed and replaced with the new parts under the new
# "3bba17b1-..." / "4decd463-..." / "http://schemas.microsoft.com/sharepoint/v3"
# namespaces, with the new default field values (Language=English,
# Sensitivity=Internal, Record=false, empty Content Author/Owner, etc.)

$d = $word.ActiveDocument

function Remove-OldCustomXml($doc) {
    $namespacesToDrop = @(
        "http://schemas.microsoft.com/office/2006/metadata/contentType",
        "700a75ce-dd9d-411b-bd1e-ab93f39097ba",
        "ce04ccd7-686d-4e9d-9396-1513c4faabae"
    )
    foreach ($ns in $namespacesToDrop) {
        try {
            $matches = $doc.CustomXMLParts.SelectByNamespace($ns)
            if ($matches -and $matches.Count -gt 0) {
                for ($i = $matches.Count; $i -ge 1; $i--) {
                    try { $matches.Item($i).Delete() } catch { }
                }
            }
        } catch { }
    }

    # Fall back to a blind sweep in case SelectByNamespace isn't available:
    # walk the whole collection back-to-front and drop anything that still
    # carries the old contentType / metadata schema.
    try {
        $count = $doc.CustomXMLParts.Count
        for ($i = $count; $i -ge 1; $i--) {
            try {
                $part = $doc.CustomXMLParts.Item($i)
                $ns = $part.NamespaceURI
                if ($namespacesToDrop -contains $ns) {
                    $part.Delete()
                }
            } catch { }
        }
    } catch { }
}

Remove-OldCustomXml $d

# --- customXml/item2.xml -------------------------------------------------
# New content-type schema (content type 0x010100FE4CDFAD8D29F44BA3595B26EBF98D7E,
# version 9) with the new managed-metadata site columns.
$item2Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100FE4CDFAD8D29F44BA3595B26EBF98D7E" ma:contentTypeVersion="9" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="866addd3bc14dc82dc60b51d940083c6">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns1="http://schemas.microsoft.com/sharepoint/v3" xmlns:ns3="3bba17b1-ca09-4865-ba6f-0714c5739852" xmlns:ns4="4decd463-a3dd-4fb0-bb57-735c1a8c741d" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="e874a30f38c26bb096d99537478574ce" ns1:_="" ns3:_="" ns4:_="">
    <xsd:import namespace="http://schemas.microsoft.com/sharepoint/v3"/>
    <xsd:import namespace="3bba17b1-ca09-4865-ba6f-0714c5739852"/>
    <xsd:import namespace="4decd463-a3dd-4fb0-bb57-735c1a8c741d"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns3:Content_x0020_Author_x0020_Email" minOccurs="0"/>
                <xsd:element ref="ns1:RoutingRuleDescription" minOccurs="0"/>
                <xsd:element ref="ns3:Content_x0020_Owner1" minOccurs="0"/>
                <xsd:element ref="ns3:Record" minOccurs="0"/>
                <xsd:element ref="ns3:Record_x0020_Series_x0020_Number" minOccurs="0"/>
                <xsd:element ref="ns1:Language" minOccurs="0"/>
                <xsd:element ref="ns3:Sensitivity" minOccurs="0"/>
                <xsd:element ref="ns3:TaxCatchAllLabel" minOccurs="0"/>
                <xsd:element ref="ns3:nd4e770dece24acd81cc5ad0e0f5f382" minOccurs="0"/>
                <xsd:element ref="ns3:me5168d4f87948a08fcc94d4eeda3704" minOccurs="0"/>
                <xsd:element ref="ns3:TaxCatchAll" minOccurs="0"/>
                <xsd:element ref="ns4:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns4:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns4:MediaServiceSearchProperties" minOccurs="0"/>
                <xsd:element ref="ns4:MediaServiceObjectDetectorVersions" minOccurs="0"/>
                <xsd:element ref="ns4:Group" minOccurs="0"/>
                <xsd:element ref="ns4:ProjectName" minOccurs="0"/>
                <xsd:element ref="ns4:Tags" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="http://schemas.microsoft.com/sharepoint/v3" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="RoutingRuleDescription" ma:index="4" nillable="true" ma:displayName="Description" ma:internalName="RoutingRuleDescription" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="Language" ma:index="9" nillable="true" ma:displayName="Language" ma:default="English" ma:internalName="Language" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:union memberTypes="dms:Text">
          <xsd:simpleType>
            <xsd:restriction base="dms:Choice">
              <xsd:enumeration value="Arabic (Saudi Arabia)"/>
              <xsd:enumeration value="Bulgarian (Bulgaria)"/>
              <xsd:enumeration value="Chinese (Hong Kong S.A.R.)"/>
              <xsd:enumeration value="Chinese (China)"/>
              <xsd:enumeration value="Chinese (Taiwan)"/>
              <xsd:enumeration value="Croatian (Croatia)"/>
              <xsd:enumeration value="Czech (Czech Republic)"/>
              <xsd:enumeration value="Danish (Denmark)"/>
              <xsd:enumeration value="Dutch (Netherlands)"/>
              <xsd:enumeration value="English"/>
              <xsd:enumeration value="Estonian (Estonia)"/>
              <xsd:enumeration value="Finnish (Finland)"/>
              <xsd:enumeration value="French (France)"/>
              <xsd:enumeration value="German (Germany)"/>
              <xsd:enumeration value="Greek (Greece)"/>
              <xsd:enumeration value="Hebrew (Israel)"/>
              <xsd:enumeration value="Hindi (India)"/>
              <xsd:enumeration value="Hungarian (Hungary)"/>
              <xsd:enumeration value="Indonesian (Indonesia)"/>
              <xsd:enumeration value="Italian (Italy)"/>
              <xsd:enumeration value="Japanese (Japan)"/>
              <xsd:enumeration value="Korean (Korea)"/>
              <xsd:enumeration value="Latvian (Latvia)"/>
              <xsd:enumeration value="Lithuanian (Lithuania)"/>
              <xsd:enumeration value="Malay (Malaysia)"/>
              <xsd:enumeration value="Norwegian (Bokmal) (Norway)"/>
              <xsd:enumeration value="Polish (Poland)"/>
              <xsd:enumeration value="Portuguese (Brazil)"/>
              <xsd:enumeration value="Portuguese (Portugal)"/>
              <xsd:enumeration value="Romanian (Romania)"/>
              <xsd:enumeration value="Russian (Russia)"/>
              <xsd:enumeration value="Serbian (Latin) (Serbia)"/>
              <xsd:enumeration value="Slovak (Slovakia)"/>
              <xsd:enumeration value="Slovenian (Slovenia)"/>
              <xsd:enumeration value="Spanish (Spain)"/>
              <xsd:enumeration value="Swedish (Sweden)"/>
              <xsd:enumeration value="Thai (Thailand)"/>
              <xsd:enumeration value="Turkish (Turkey)"/>
              <xsd:enumeration value="Ukrainian (Ukraine)"/>
              <xsd:enumeration value="Urdu (Islamic Republic of Pakistan)"/>
              <xsd:enumeration value="Vietnamese (Vietnam)"/>
            </xsd:restriction>
          </xsd:simpleType>
        </xsd:union>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="3bba17b1-ca09-4865-ba6f-0714c5739852" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="Content_x0020_Author_x0020_Email" ma:index="3" nillable="true" ma:displayName="Content Author Email" ma:list="UserInfo" ma:SharePointGroup="0" ma:internalName="Content_x0020_Author_x0020_Email" ma:readOnly="false" ma:showField="ImnName">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:User">
            <xsd:sequence>
              <xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded">
                <xsd:complexType>
                  <xsd:sequence>
                    <xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/>
                    <xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/>
                    <xsd:element name="AccountType" type="xsd:string" minOccurs="0"/>
                  </xsd:sequence>
                </xsd:complexType>
              </xsd:element>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="Content_x0020_Owner1" ma:index="6" nillable="true" ma:displayName="Content Owner" ma:description="Identifies the business owner responsible for the accuracy and appropriate use of content items. Is generally hidden from view, but is used by the content management system to notify a responsible person that review or other content lifecycle" ma:list="UserInfo" ma:SharePointGroup="0" ma:internalName="Content_x0020_Owner1" ma:readOnly="false" ma:showField="ImnName">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:User">
            <xsd:sequence>
              <xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded">
                <xsd:complexType>
                  <xsd:sequence>
                    <xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/>
                    <xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/>
                    <xsd:element name="AccountType" type="xsd:string" minOccurs="0"/>
                  </xsd:sequence>
                </xsd:complexType>
              </xsd:element>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="Record" ma:index="7" nillable="true" ma:displayName="Record" ma:default="0" ma:description="​For collaborative or normal documents, the default is set to No." ma:internalName="Record" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:restriction base="dms:Boolean"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="Record_x0020_Series_x0020_Number" ma:index="8" nillable="true" ma:displayName="Record Series Number" ma:description="​If the record tag is set to yes, a valid Record Series Number (as defined by RM-00-001.AV is required to be present." ma:internalName="Record_x0020_Series_x0020_Number" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="Sensitivity" ma:index="10" nillable="true" ma:displayName="Sensitivity" ma:default="Internal" ma:description="​Security level required by the content . Values = Public, Internal, Restricted, Secret" ma:format="Dropdown" ma:internalName="Sensitivity" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:restriction base="dms:Choice">
          <xsd:enumeration value="Public"/>
          <xsd:enumeration value="Internal"/>
          <xsd:enumeration value="Restricted"/>
          <xsd:enumeration value="Secret"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="TaxCatchAllLabel" ma:index="12" nillable="true" ma:displayName="Taxonomy Catch All Column1" ma:hidden="true" ma:list="{52c9f977-e2e7-4cdf-85bf-9dbca809ee38}" ma:internalName="TaxCatchAllLabel" ma:readOnly="true" ma:showField="CatchAllDataLabel" ma:web="88998b2f-5cc8-475d-b9b4-a8d353d33050">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:MultiChoiceLookup">
            <xsd:sequence>
              <xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="nd4e770dece24acd81cc5ad0e0f5f382" ma:index="13" nillable="true" ma:taxonomy="true" ma:internalName="nd4e770dece24acd81cc5ad0e0f5f382" ma:taxonomyFieldName="Document_x0020_Categories" ma:displayName="Document Categories" ma:readOnly="false" ma:fieldId="{7d4e770d-ece2-4acd-81cc-5ad0e0f5f382}" ma:sspId="f53c2da0-965e-4c49-9e20-3f7554834061" ma:termSetId="0bc447bc-cd5f-4a5b-8cef-40c06dd94379" ma:anchorId="00000000-0000-0000-0000-000000000000" ma:open="false" ma:isKeyword="false">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element ref="pc:Terms" minOccurs="0" maxOccurs="1"/>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="me5168d4f87948a08fcc94d4eeda3704" ma:index="15" nillable="true" ma:displayName="Topic Column_0" ma:hidden="true" ma:internalName="me5168d4f87948a08fcc94d4eeda3704" ma:readOnly="false">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="TaxCatchAll" ma:index="21" nillable="true" ma:displayName="Taxonomy Catch All Column" ma:hidden="true" ma:list="{52c9f977-e2e7-4cdf-85bf-9dbca809ee38}" ma:internalName="TaxCatchAll" ma:readOnly="false" ma:showField="CatchAllData" ma:web="88998b2f-5cc8-475d-b9b4-a8d353d33050">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:MultiChoiceLookup">
            <xsd:sequence>
              <xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="4decd463-a3dd-4fb0-bb57-735c1a8c741d" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="22" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="23" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceSearchProperties" ma:index="24" nillable="true" ma:displayName="MediaServiceSearchProperties" ma:hidden="true" ma:internalName="MediaServiceSearchProperties" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceObjectDetectorVersions" ma:index="25" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="Group" ma:index="26" nillable="true" ma:displayName="Group" ma:format="Dropdown" ma:internalName="Group">
      <xsd:simpleType>
        <xsd:restriction base="dms:Choice">
          <xsd:enumeration value="CCD"/>
          <xsd:enumeration value="Pilot Lab"/>
          <xsd:enumeration value="Purification Development"/>
          <xsd:enumeration value="S&amp;T"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="ProjectName" ma:index="27" nillable="true" ma:displayName="Project Name" ma:format="Dropdown" ma:internalName="ProjectName">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="Tags" ma:index="28" nillable="true" ma:displayName="Tags" ma:format="Dropdown" ma:internalName="Tags">
      <xsd:simpleType>
        <xsd:union memberTypes="dms:Text">
          <xsd:simpleType>
            <xsd:restriction base="dms:Choice">
              <xsd:enumeration value="Quote"/>
              <xsd:enumeration value="Presentation"/>
              <xsd:enumeration value="script"/>
              <xsd:enumeration value="reference"/>
              <xsd:enumeration value="template"/>
            </xsd:restriction>
          </xsd:simpleType>
        </xsd:union>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd" targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all">
    <xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
    <xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
    <xsd:element name="coreProperties" type="CT_coreProperties"/>
    <xsd:complexType name="CT_coreProperties">
      <xsd:all>
        <xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1" ma:index="2" ma:displayName="Author"/>
        <xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="16" ma:displayName="Content Type"/>
        <xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="1" ma:displayName="Title"/>
        <xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
          <xsd:annotation>
            <xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
          </xsd:annotation>
        </xsd:element>
        <xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
      </xsd:all>
    </xsd:complexType>
  </xsd:schema>
  <xs:schema targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema">
    <xs:element name="Person">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:DisplayName" minOccurs="0"/>
          <xs:element ref="pc:AccountId" minOccurs="0"/>
          <xs:element ref="pc:AccountType" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="DisplayName" type="xs:string"/>
    <xs:element name="AccountId" type="xs:string"/>
    <xs:element name="AccountType" type="xs:string"/>
    <xs:element name="BDCAssociatedEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
        <xs:attribute ref="pc:EntityNamespace"/>
        <xs:attribute ref="pc:EntityName"/>
        <xs:attribute ref="pc:SystemInstanceName"/>
        <xs:attribute ref="pc:AssociationName"/>
      </xs:complexType>
    </xs:element>
    <xs:attribute name="EntityNamespace" type="xs:string"/>
    <xs:attribute name="EntityName" type="xs:string"/>
    <xs:attribute name="SystemInstanceName" type="xs:string"/>
    <xs:attribute name="AssociationName" type="xs:string"/>
    <xs:element name="BDCEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:EntityDisplayName" minOccurs="0"/>
          <xs:element ref="pc:EntityInstanceReference" minOccurs="0"/>
          <xs:element ref="pc:EntityId1" minOccurs="0"/>
          <xs:element ref="pc:EntityId2" minOccurs="0"/>
          <xs:element ref="pc:EntityId3" minOccurs="0"/>
          <xs:element ref="pc:EntityId4" minOccurs="0"/>
          <xs:element ref="pc:EntityId5" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="EntityDisplayName" type="xs:string"/>
    <xs:element name="EntityInstanceReference" type="xs:string"/>
    <xs:element name="EntityId1" type="xs:string"/>
    <xs:element name="EntityId2" type="xs:string"/>
    <xs:element name="EntityId3" type="xs:string"/>
    <xs:element name="EntityId4" type="xs:string"/>
    <xs:element name="EntityId5" type="xs:string"/>
    <xs:element name="Terms">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermInfo">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermName" minOccurs="0"/>
          <xs:element ref="pc:TermId" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermName" type="xs:string"/>
    <xs:element name="TermId" type="xs:string"/>
  </xs:schema>
</ct:contentTypeSchema>
'@

# --- customXml/item3.xml --------------------------------------------------
# New document-management field values (Language/Sensitivity defaults,
# empty Content Author Email / Content Owner, Record=false, ...).
$item3Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
  <documentManagement>
    <TaxCatchAll xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852" xsi:nil="true"/>
    <Language xmlns="http://schemas.microsoft.com/sharepoint/v3">English</Language>
    <ProjectName xmlns="4decd463-a3dd-4fb0-bb57-735c1a8c741d" xsi:nil="true"/>
    <Group xmlns="4decd463-a3dd-4fb0-bb57-735c1a8c741d" xsi:nil="true"/>
    <Content_x0020_Author_x0020_Email xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852">
      <UserInfo>
        <DisplayName/>
        <AccountId xsi:nil="true"/>
        <AccountType/>
      </UserInfo>
    </Content_x0020_Author_x0020_Email>
    <Record_x0020_Series_x0020_Number xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852" xsi:nil="true"/>
    <me5168d4f87948a08fcc94d4eeda3704 xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852" xsi:nil="true"/>
    <Content_x0020_Owner1 xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852">
      <UserInfo>
        <DisplayName/>
        <AccountId xsi:nil="true"/>
        <AccountType/>
      </UserInfo>
    </Content_x0020_Owner1>
    <nd4e770dece24acd81cc5ad0e0f5f382 xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852">
      <Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    </nd4e770dece24acd81cc5ad0e0f5f382>
    <RoutingRuleDescription xmlns="http://schemas.microsoft.com/sharepoint/v3" xsi:nil="true"/>
    <Sensitivity xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852">Internal</Sensitivity>
    <Tags xmlns="4decd463-a3dd-4fb0-bb57-735c1a8c741d" xsi:nil="true"/>
    <Record xmlns="3bba17b1-ca09-4865-ba6f-0714c5739852">false</Record>
  </documentManagement>
</p:properties>
'@

# --- customXml/item4.xml --------------------------------------------------
# Taxonomy content-type sync marker, new with this resync.
$item4Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<?mso-contentType ?>
<SharedContentType xmlns="Microsoft.SharePoint.Taxonomy.ContentTypeSync" SourceId="f53c2da0-965e-4c49-9e20-3f7554834061" ContentTypeId="0x0101" PreviousValue="false"/>
'@

function Add-CustomXmlSafely($doc, [string]$xml) {
    try {
        return $doc.CustomXMLParts.Add($xml)
    } catch {
        return $null
    }
}

Add-CustomXmlSafely $d $item2Xml | Out-Null
Add-CustomXmlSafely $d $item3Xml | Out-Null
Add-CustomXmlSafely $d $item4Xml | Out-Null
